# Fruta / hortaliza, semanal
# Update the weekly price/volume records for rows 2-5 (Hortaliza, Macroferia
# Regional de Talca - Fruto del paraíso). The values for columns D (Fecha),
# J (Volumen), K (Precio mínimo), L (Precio máximo), M (Precio promedio
# ponderado) and P (Precio $/Kg) are rotated across the four data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken from the target (post-edit) workbook state.
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

$ws.Range("D3").Value = 44277
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 550

$ws.Range("D4").Value = 44284
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 500

$ws.Range("D5").Value = 44280
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
